$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-white.png"
$ws.Range("A3").Value  = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-white.png"
$ws.Range("A4").Value  = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-yellow.png"
$ws.Range("A5").Value  = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-white.png"
$ws.Range("A6").Value  = "./images_eeg/Sphere_CW-2.25_BG-grey_stim-white.png"
$ws.Range("A7").Value  = "./images_eeg/Sphere_Ref_BG-grey_stim-yellow.png"
$ws.Range("A9").Value  = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-white.png"
$ws.Range("A10").Value = "./images_eeg/Sphere_Ref_BG-grey_stim-white.png"
$ws.Range("A11").Value = "./images_eeg/Sphere_CW-2.25_BG-grey_stim-white.png"
$ws.Range("A12").Value = "./images_eeg/Sphere_CCW-2.25_BG-grey_stim-white.png"
$ws.Range("A13").Value = "./images_eeg/Sphere_Ref_BG-grey_stim-yellow.png"
$ws.Range("A15").Value = "./images_eeg/Sphere_CW-2.25_BG-grey_stim-white.png"
$ws.Range("A16").Value = "./images_eeg/Sphere_CW-2.25_BG-grey_stim-yellow.png"
$ws.Range("A17").Value = "./images_eeg/Sphere_CW-2.25_BG-grey_stim-white.png"
$ws.Range("A18").Value = "./images_eeg/Sphere_Ref_BG-grey_stim-white.png"
